$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 0
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 1
$ws.Range("H5").Value = 0
